$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124, shifting existing rows 124:181 down to 125:182
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new record's data
$ws.Range("A124").Value = 8
$ws.Range("B124").Value = "Terminal La Palmera de La Serena"
$ws.Range("C124").Value = "Coquimbo"
$ws.Range("D124").Value = 44596
$ws.Range("E124").Value = 4
$ws.Range("F124").Value = 100112031
$ws.Range("G124").Value = "Poroto verde"
$ws.Range("H124").Value = "Magnum"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 400
$ws.Range("K124").Value = 37000
$ws.Range("L124").Value = 38000
$ws.Range("M124").Value = 37500
$ws.Range("N124").Value = "$/malla 25 kilos"
$ws.Range("O124").Value = "Provincia de Limarí"
$ws.Range("P124").Value = 1500
$ws.Range("Q124").Value = 25
$ws.Range("R124").Value = "Hortaliza"
